$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.8626097691870205
$ws.Range("D2").Value = 0.3976581394938257

$ws.Range("C3").Value = 0.4396111220284893
$ws.Range("D3").Value = 0.6645069094021108

$ws.Range("C4").Value = 0.8628456531644179
$ws.Range("D4").Value = 0.3975312190896196

$ws.Range("C5").Value = -0.4660211551627736
$ws.Range("D5").Value = 0.6457814353122195

$ws.Range("C6").Value = 1.079854876138309
$ws.Range("D6").Value = 0.2919036008572575

$ws.Range("C7").Value = 1.639722369904767
$ws.Range("D7").Value = 0.1152852269545397

$ws.Range("C8").Value = 0.7042833850380079
$ws.Range("D8").Value = 0.4886450673384424

$ws.Range("C9").Value = 0.6097881176744524
$ws.Range("D9").Value = 0.5482493192642757

$ws.Range("C10").Value = -0.7293232672179952
$ws.Range("D10").Value = 0.4734978824273424

$ws.Range("C11").Value = -1.094718551873542
$ws.Range("D11").Value = 0.2854808967243376
